$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.196023225784302
$ws.Range("B1").Value = 2.481930017471313
$ws.Range("C1").Value = 4.074303150177002
$ws.Range("D1").Value = 2.103570938110352
$ws.Range("E1").Value = 1.1844722032547
